$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 38 and 39 swap coin identity (name/link) and get new price/volume values.
# NumberFormat is forced to Text ("@") on each Price cell before the write so the
# numeric-looking strings (e.g. "2.12", "47.32") are kept as literal text instead
# of being coerced into numbers (which would drop meaningful trailing zeros).
$ws.Range("D38").NumberFormat = "@"
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").Value = "2.12"
$ws.Range("E38").Value = "  +3.20%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("B39").Value = "Arweave"
$ws.Range("C39").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D39").Value = "47.32"
$ws.Range("E39").Value = "  -2.66%  "

# Remaining price/volume updates
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.006.05"
$ws.Range("E2").Value = "  +1.19%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.116.09"
$ws.Range("E3").Value = "  +2.17%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.00"
$ws.Range("E5").Value = "  +0.60%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.96"
$ws.Range("E6").Value = "  +2.19%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.113.20"
$ws.Range("E8").Value = "  +2.41%  "
$ws.Range("E9").Value = "  +0.47%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.46"
$ws.Range("E10").Value = "  -3.26%  "
$ws.Range("E11").Value = "  +1.38%  "
$ws.Range("E12").Value = "  +0.30%  "
$ws.Range("E13").Value = "  +0.78%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.35"
$ws.Range("E14").Value = "  +0.51%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.125"
$ws.Range("E15").Value = "  +0.26%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.630.88"
$ws.Range("E16").Value = "  +2.14%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.066.42"
$ws.Range("E17").Value = "  +1.21%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.21"
$ws.Range("E18").Value = "  -0.15%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.114.74"
$ws.Range("E19").Value = "  +2.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.36"
$ws.Range("E20").Value = "  +0.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "486.89"
$ws.Range("E21").Value = "  +4.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.721"
$ws.Range("E22").Value = "  +1.96%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.58"
$ws.Range("E23").Value = "  +1.23%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.60"
$ws.Range("E24").Value = "  +1.47%  "
$ws.Range("E25").Value = "  +2.68%  "
$ws.Range("E26").Value = "  +3.74%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.10"
$ws.Range("E27").Value = "  -0.17%  "
$ws.Range("E28").Value = "  -0.21%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.06"
$ws.Range("E29").Value = "  -3.94%  "
$ws.Range("E30").Value = "  -1.80%  "
$ws.Range("E31").Value = "  +2.10%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "29.09"
$ws.Range("E32").Value = "  +2.79%  "
$ws.Range("E33").Value = "  +0.46%  "
$ws.Range("E34").Value = "  -3.52%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.97"
$ws.Range("E36").Value = "  +2.08%  "
$ws.Range("E37").Value = "  -0.38%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.16"
$ws.Range("E40").Value = "  +1.31%  "
$ws.Range("E41").Value = "  +1.89%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.123"
$ws.Range("E42").Value = "  +1.58%  "
$ws.Range("E43").Value = "  +0.68%  "
$ws.Range("E44").Value = "  -2.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.844.34"
$ws.Range("E45").Value = "  +3.58%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "386.81"
$ws.Range("E46").Value = "  +0.49%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0360"
$ws.Range("E47").Value = "  +0.10%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "136.32"
$ws.Range("E48").Value = "  +1.40%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.11"
$ws.Range("E50").Value = "  +0.51%  "
$ws.Range("E51").Value = "  -0.21%  "
